# Rename the column headers in row 1 of the "AHB-Diff" sheet so that the
# former "_old"/"_new" suffixes are replaced by the respective format
# version suffixes "_FV2310" / "_FV2404" (the "diff" column is untouched),
# freeze the header row, and wrap the data range in an Excel Table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Freeze the header row (split below row 1, freeze panes on).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into an Excel Table ("Table1") with an autofilter,
# matching the worksheet dimension A1:U68.
$tblRange = $ws.Range("A1:U68")
$tbl = $ws.ListObjects.Add(1, $tblRange, [System.Type]::Missing, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

Write-Host "Header renaming, freeze panes, and table creation complete."
